# Reference-date correction: change "2014"/"2015" citation years to "2021"
# on three slides, without disturbing existing run/line-break structure.
#
# Note: TextRange.Replace() in this runtime does not cleanly delete the
# matched text before inserting the replacement, so instead we locate the
# exact substring via TextRange.Text / String.IndexOf, then grab that exact
# span with TextRange.Characters(start, length) and overwrite its .Text.
# When the span lines up with an existing run's full text, this updates the
# run in place instead of fragmenting it.
# (Named parameters aren't supported by this host's PowerShell subset, so
# the helper uses positional parameters only.)

function Set-SubstringText($TextRange, $OldText, $NewText) {
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Substring not found: $OldText"
    }
    $sub = $TextRange.Characters($idx + 1, $OldText.Length)
    $sub.Text = $NewText
}

$p = $ppt.ActivePresentation

# Slide 2 (sldId 258): title "What is a Just, Learning Culture (KIM, Et Al. 2014)"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
Set-SubstringText $tr2 "(KIM, Et Al. 2014)" "(KIM, Et Al. 2021)"

# Slide 4 (sldId 260): title "Bad Apple Theory (kim, et al. 2014)"
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(1).TextFrame.TextRange
Set-SubstringText $tr4 ", et al. 2014)" ", et al. 2021)"

# Slide 11 (sldId 257): references placeholder, DevOps Handbook citation year
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item(2).TextFrame.TextRange
$oldRef = ", P., Willis, J., & Humble, J. (2015). The DevOps Handbook: How to Create World-Class Speed, Reliability, and Security in Technology Organizations. It Revolution Press."
$newRef = ", P., Willis, J., & Humble, J. (2021). The DevOps Handbook: How to Create World-Class Speed, Reliability, and Security in Technology Organizations. It Revolution Press."
Set-SubstringText $tr11 $oldRef $newRef
